$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.566.16'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '2.443.63'
$ws.Range("E3").Value = '  -2.30%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '547.46'
$ws.Range("E5").Value = '  -0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.89'
$ws.Range("E6").Value = '  -0.99%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.583'
$ws.Range("E8").Value = '  -1.90%  '
$ws.Range("D9").Value = '2.442.28'
$ws.Range("E9").Value = '  -2.31%  '
$ws.Range("E10").Value = '  -0.98%  '
$ws.Range("E11").Value = '  +1.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.45'
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("E13").Value = '  -2.12%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.87'
$ws.Range("E14").Value = '  -1.09%  '
$ws.Range("D15").Value = '2.881.04'
$ws.Range("E15").Value = '  -2.11%  '
$ws.Range("E16").Value = '  +2.81%  '
$ws.Range("D17").Value = '61.369.00'
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("D18").Value = '2.439.47'
$ws.Range("E18").Value = '  -1.64%  '
$ws.Range("E19").Value = '  -3.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.89'
$ws.Range("E20").Value = '  -1.81%  '
$ws.Range("E21").Value = '  -1.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '318.34'
$ws.Range("E22").Value = '  -0.87%  '
$ws.Range("E23").Value = '  -0.12%  '
$ws.Range("E24").Value = '  +5.26%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.87'
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("D26").Value = '0.0₃0971'
$ws.Range("E26").Value = '  -5.65%  '
$ws.Range("D27").Value = '2.561.76'
$ws.Range("E27").Value = '  -1.71%  '
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("E29").Value = '  -1.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.81'
$ws.Range("E30").Value = '  +2.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '524.98'
$ws.Range("E31").Value = '  -2.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.17'
$ws.Range("E32").Value = '  -2.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.147'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("E35").Value = '  +0.93%  '
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.62'
$ws.Range("E37").Value = '  -4.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.73'
$ws.Range("E38").Value = '  -2.71%  '
$ws.Range("E39").Value = '  +1.19%  '
$ws.Range("E40").Value = '  -1.63%  '
$ws.Range("E41").Value = '  +2.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '138.65'
$ws.Range("E42").Value = '  -4.03%  '
$ws.Range("E43").Value = '  +0.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.35'
$ws.Range("E44").Value = '  -0.49%  '
$ws.Range("E45").Value = '  -1.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '142.52'
$ws.Range("E46").Value = '  -4.12%  '
$ws.Range("E47").Value = '  +0.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.16'
$ws.Range("E48").Value = '  +0.63%  '
$ws.Range("E49").Value = '  -1.80%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.587'
$ws.Range("E50").Value = '  -0.75%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0929'
$ws.Range("E51").Value = '  -1.53%  '
